$d = $word.ActiveDocument

# Remove the existing _GoBack bookmark (it currently sits, collapsed,
# at the very end of the "Testing testing" paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Split the last paragraph ("Testing testing") so that a brand new
# paragraph containing "Test 1234" follows it.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$splitPoint = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)
$splitPoint.InsertAfter("`r")

# Fill in the new paragraph's text.
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.InsertBefore("Test 1234")

# Re-create the _GoBack bookmark spanning from the very start of the
# document to the very end (after "Test 1234"), so bookmarkStart lands
# right before the first run and bookmarkEnd lands at the end of the
# document.
$fullRange = $d.Range(0, $d.Content.End)
$d.Bookmarks.Add("_GoBack", $fullRange)
